$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.741.68'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.634.93'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  -0.15%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '215.22'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  -0.13%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.0635'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.12%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '19.54'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -4.03%  '
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.636.75'
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.859.91'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("D16").Value = '0.0₃0764'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").Value = '25.764.78'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("E20").Value = '  +1.30%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '193.95'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("E22").Value = '  +0.03%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.27'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +2.19%  '
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("E25").Value = '  +2.60%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '140.39'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.49%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.121'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("E28").Value = '  +1.05%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.51'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  -0.09%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.0492'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.36%  '
$ws.Range("E32").Value = '  +1.29%  '
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("D37").Value = '1.122.17'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.547'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.35%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.52'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.86%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0155'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.07%  '
$ws.Range("E41").Value = '  +0.56%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.57'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.97%  '
$ws.Range("E43").Value = '  +0.63%  '
$ws.Range("E44").Value = '  +0.14%  '
$ws.Range("D45").Value = '1.769.10'
$ws.Range("E45").Value = '  -0.50%  '
$ws.Range("E46").Value = '  -2.18%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '55.02'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("E48").Value = '  -2.29%  '
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("E50").Value = '  -3.16%  '
$ws.Range("E51").Value = '  +0.87%  '
